$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that sit between "Voltage Regulator 7.4V-5V" (old row 12)
# and the blank spacer (old row 13) and the "Ordered" section (old row 16).
# This shifts rows 16-21 up to 14-19 and auto-adjusts the D20/D21 formulas.
$ws.Rows("12:13").Delete()

# Row 9: "Things to ask about:" -> "Received from University:" and clear its price
$ws.Range("C9").Value = "Received from University:"
$ws.Range("D9").ClearContents()

# Row 10: ON/OFF Switch price reset to 0
$ws.Range("D10").Value = 0

# Row 11: was "ADS1115 ADC 16-bit ADC" / 2.69 -> now "Voltage Regulator 7.4V-5V" / 0
$ws.Range("C11").Value = "Voltage Regulator 7.4V-5V"
$ws.Range("D11").Value = 0

# The Geophone SM-24 hyperlink used to live on E17; after the row deletion its
# text moved to E15, but the hyperlink anchor needs to be re-pointed there too.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E15"), "https://www.aliexpress.com/item/1005005275736468.html?spm=a2g0o.productlist.main.1.3a1072cewm2hWM&algo_pvid=b268c8e1-af98-4201-8858-09cd4e6887e0&algo_exp_id=b268c8e1-af98-4201-8858-09cd4e6887e0-0&pdp_npi=4%40dis%21CAD%2173.14%2171.68%21%21%2152.61%2151.56%21%40210307bf17265994460883461ec49b%2112000038407879035%21sea%21CA%210%21ABX&curPageLogUid=dF82d1Pl71AG&utparam-url=scene%3Asearch%7Cquery_from%3A", "", "", "https://www.aliexpress.com/item/1005005275736468.html?spm=a2g0o.productlist.main.1.3a1072cewm2hWM&algo_pvid=b268c8e1-af98-4201-8858-09cd4e6887e0&algo_exp_id=b268c8e1-af98-4201-8858-09cd4e6887e0-0&pdp_npi=4%40dis%21CAD%2173.14%2171.68%21%21%2152.61%2151.56%21%40210307bf17265994460883461ec49b%2112000038407879035%21sea%21CA%210%21ABX&curPageLogUid=dF82d1Pl71AG&utparam-url=scene%3Asearch%7Cquery_from%3A")
$ws.Range("E15").Style = "Hyperlink"

# Update the selected cell shown in the saved workbook
$ws.Range("D28").Select()
